# refatoracao - calculos de apoio medio
#
# Rebuilds the descriptive-stats sheet with the new "apoio_*" and
# "*_std/min/max" breakdown columns:
#   media_sucesso/std_sucesso/min_sucesso/max_sucesso -> arrecadado_avg/std/min/max
#   apoio_medio (values recomputed) + new apoio_std/min/max
#   media_contribuicoes -> contribuicoes_med + new contribuicoes_std/min/max
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1) Header row (A1:V1) ----
$headers = New-Object 'object[,]' 1,22
$headers[0,0] = "modalidade"
$headers[0,1] = "mencao"
$headers[0,2] = "total"
$headers[0,3] = "total_sucesso"
$headers[0,4] = "particip"
$headers[0,5] = "taxa_sucesso"
$headers[0,6] = "arrecadado_sucesso"
$headers[0,7] = "arrecadado_avg"
$headers[0,8] = "arrecadado_std"
$headers[0,9] = "arrecadado_min"
$headers[0,10] = "arrecadado_max"
$headers[0,11] = "apoio_medio"
$headers[0,12] = "apoio_std"
$headers[0,13] = "apoio_min"
$headers[0,14] = "apoio_max"
$headers[0,15] = "contribuicoes"
$headers[0,16] = "contribuicoes_med"
$headers[0,17] = "contribuicoes_std"
$headers[0,18] = "contribuicoes_min"
$headers[0,19] = "contribuicoes_max"
$headers[0,20] = "menor_ano"
$headers[0,21] = "maior_ano"
$ws.Range("A1:V1").Value = $headers

# Copy the existing header style (bold/border/centered, s=4 from L1) onto
# the newly introduced header cells Q1:V1.
$ws.Range("L1").Copy()
$ws.Range("Q1:V1").PasteSpecial(-4122)

# ---- 2) Data rows (A2:V22) ----
$data = New-Object 'object[,]' 21,22
$data[0,0] = "aon"
$data[0,1] = "angelo_agostini"
$data[0,2] = 67
$data[0,3] = 61
$data[0,4] = 0.050187265917603
$data[0,5] = 0.9104477611940298
$data[0,6] = 3045252.806306844
$data[0,7] = 49922.17715257122
$data[0,8] = 101183.4689374139
$data[0,9] = 2944.086470789134
$data[0,10] = 679297.6600721752
$data[0,11] = 78.57672143399174
$data[0,12] = 23.63877229505447
$data[0,13] = 32.2694534583262
$data[0,14] = 151.1292159501072
$data[0,15] = 35293
$data[0,16] = 578.5737704918033
$data[0,17] = 929.9762624034444
$data[0,18] = 55
$data[0,19] = 6494
$data[0,20] = 2013
$data[0,21] = 2023
$data[1,0] = "aon"
$data[1,1] = "ccxp"
$data[1,2] = 156
$data[1,3] = 126
$data[1,4] = 0.1168539325842697
$data[1,5] = 0.8076923076923077
$data[1,6] = 2860818.438596986
$data[1,7] = 22704.90824283323
$data[1,8] = 20582.18081328499
$data[1,9] = 1720.659275370021
$data[1,10] = 154365.9837040891
$data[1,11] = 87.11000777143747
$data[1,12] = 35.07795799700576
$data[1,13] = 33.80063482849972
$data[1,14] = 257.7853211115706
$data[1,15] = 34535
$data[1,16] = 274.0873015873016
$data[1,17] = 239.4804716828918
$data[1,18] = 36
$data[1,19] = 1815
$data[1,20] = 2014
$data[1,21] = 2023
$data[2,0] = "aon"
$data[2,1] = "disputa"
$data[2,2] = 400
$data[2,3] = 245
$data[2,4] = 0.299625468164794
$data[2,5] = 0.6125
$data[2,6] = 6532941.879185004
$data[2,7] = 26665.06889463267
$data[2,8] = 27722.17615151563
$data[2,9] = 787.1021772339901
$data[2,10] = 163173.270269744
$data[2,11] = 95.92481440598637
$data[2,12] = 62.28205825061158
$data[2,13] = 21.61624650544615
$data[2,14] = 792.0360759681182
$data[2,15] = 70527
$data[2,16] = 287.865306122449
$data[2,17] = 263.6097661931087
$data[2,18] = 11
$data[2,19] = 1588
$data[2,20] = 2012
$data[2,21] = 2023
$data[3,0] = "aon"
$data[3,1] = "erotismo"
$data[3,2] = 123
$data[3,3] = 82
$data[3,4] = 0.09213483146067415
$data[3,5] = 0.6666666666666666
$data[3,6] = 2306077.727819387
$data[3,7] = 28122.89911974862
$data[3,8] = 23954.88805282878
$data[3,9] = 1990.9450708267
$data[3,10] = 125535.7448333134
$data[3,11] = 94.90136016654937
$data[3,12] = 42.06922775736284
$data[3,13] = 37.30506273801686
$data[3,14] = 323.2845357010965
$data[3,15] = 25810
$data[3,16] = 314.7560975609756
$data[3,17] = 268.2786261681259
$data[3,18] = 22
$data[3,19] = 1539
$data[3,20] = 2012
$data[3,21] = 2023
$data[4,0] = "aon"
$data[4,1] = "fantasia"
$data[4,2] = 279
$data[4,3] = 180
$data[4,4] = 0.2089887640449438
$data[4,5] = 0.6451612903225806
$data[4,6] = 4816835.598395908
$data[4,7] = 26760.19776886616
$data[4,8] = 30295.08564925056
$data[4,9] = 94.898114598278
$data[4,10] = 264585.9073482947
$data[4,11] = 88.02163529111277
$data[4,12] = 33.79527781109714
$data[4,13] = 23.15006403629383
$data[4,14] = 213.9734252018395
$data[4,15] = 54018
$data[4,16] = 300.1
$data[4,17] = 283.7555817382237
$data[4,18] = 3
$data[4,19] = 1711
$data[4,20] = 2012
$data[4,21] = 2023
$data[5,0] = "aon"
$data[5,1] = "ficcao_cientifica"
$data[5,2] = 296
$data[5,3] = 179
$data[5,4] = 0.2217228464419476
$data[5,5] = 0.6047297297297297
$data[5,6] = 5401646.146659081
$data[5,7] = 30176.79411541386
$data[5,8] = 46623.43714083682
$data[5,9] = 54.53892516702949
$data[5,10] = 537544.5528256212
$data[5,11] = 94.68234986858575
$data[5,12] = 69.67419713686293
$data[5,13] = 30.56837093393595
$data[5,14] = 792.0360759681182
$data[5,15] = 61529
$data[5,16] = 343.7374301675978
$data[5,17] = 493.6165473708339
$data[5,18] = 1
$data[5,19] = 5879
$data[5,20] = 2012
$data[5,21] = 2023
$data[6,0] = "aon"
$data[6,1] = "fiq"
$data[6,2] = 219
$data[6,3] = 162
$data[6,4] = 0.1640449438202247
$data[6,5] = 0.7397260273972602
$data[6,6] = 4821257.849740589
$data[6,7] = 29760.85092432462
$data[6,8] = 38181.97911034847
$data[6,9] = 1405.425637867093
$data[6,10] = 396557.4961875453
$data[6,11] = 86.76189291082824
$data[6,12] = 32.20538078565799
$data[6,13] = 21.61624650544615
$data[6,14] = 199.8601709743299
$data[6,15] = 54948
$data[6,16] = 339.1851851851852
$data[6,17] = 351.9795170180732
$data[6,18] = 31
$data[6,19] = 3266
$data[6,20] = 2011
$data[6,21] = 2023
$data[7,0] = "aon"
$data[7,1] = "folclore"
$data[7,2] = 200
$data[7,3] = 140
$data[7,4] = 0.149812734082397
$data[7,5] = 0.7
$data[7,6] = 4507319.894909304
$data[7,7] = 32195.14210649503
$data[7,8] = 45455.39030881756
$data[7,9] = 1081.472570086762
$data[7,10] = 396557.4961875453
$data[7,11] = 90.34545930449829
$data[7,12] = 67.62639234759273
$data[7,13] = 32.2694534583262
$data[7,14] = 792.0360759681182
$data[7,15] = 51967
$data[7,16] = 371.1928571428571
$data[7,17] = 416.4698145405411
$data[7,18] = 17
$data[7,19] = 3266
$data[7,20] = 2012
$data[7,21] = 2023
$data[8,0] = "aon"
$data[8,1] = "herois"
$data[8,2] = 276
$data[8,3] = 157
$data[8,4] = 0.2067415730337079
$data[8,5] = 0.5688405797101449
$data[8,6] = 3474485.413264631
$data[8,7] = 22130.48033926517
$data[8,8] = 25554.86463253136
$data[8,9] = 989.6825900594964
$data[8,10] = 161153.6281071717
$data[8,11] = 96.01288922432792
$data[8,12] = 67.13300946444733
$data[8,13] = 21.61624650544615
$data[8,14] = 792.0360759681182
$data[8,15] = 38102
$data[8,16] = 242.687898089172
$data[8,17] = 251.542596013111
$data[8,18] = 11
$data[8,19] = 1588
$data[8,20] = 2012
$data[8,21] = 2023
$data[9,0] = "aon"
$data[9,1] = "hqmix"
$data[9,2] = 122
$data[9,3] = 108
$data[9,4] = 0.09138576779026217
$data[9,5] = 0.8852459016393442
$data[9,6] = 3759007.866591265
$data[9,7] = 34805.62839436356
$data[9,8] = 76815.19294381328
$data[9,9] = 787.1021772339901
$data[9,10] = 679297.6600721752
$data[9,11] = 81.01947145980097
$data[9,12] = 26.64100216145226
$data[9,13] = 21.61624650544615
$data[9,14] = 172.1642729447236
$data[9,15] = 42155
$data[9,16] = 390.3240740740741
$data[9,17] = 706.124751741857
$data[9,18] = 12
$data[9,19] = 6494
$data[9,20] = 2013
$data[9,21] = 2023
$data[10,0] = "aon"
$data[10,1] = "hqmix"
$data[10,2] = 122
$data[10,3] = 108
$data[10,4] = 0.09138576779026217
$data[10,5] = 0.8852459016393442
$data[10,6] = 3759007.866591265
$data[10,7] = 34805.62839436356
$data[10,8] = 76815.19294381328
$data[10,9] = 787.1021772339901
$data[10,10] = 679297.6600721752
$data[10,11] = 81.01947145980097
$data[10,12] = 26.64100216145226
$data[10,13] = 21.61624650544615
$data[10,14] = 172.1642729447236
$data[10,15] = 42155
$data[10,16] = 390.3240740740741
$data[10,17] = 706.124751741857
$data[10,18] = 12
$data[10,19] = 6494
$data[10,20] = 2013
$data[10,21] = 2023
$data[11,0] = "aon"
$data[11,1] = "jogos"
$data[11,2] = 284
$data[11,3] = 201
$data[11,4] = 0.2127340823970038
$data[11,5] = 0.7077464788732394
$data[11,6] = 5604766.32132519
$data[11,7] = 27884.40955883179
$data[11,8] = 51240.25465965582
$data[11,9] = 1411.863595952828
$data[11,10] = 679297.6600721752
$data[11,11] = 91.15737775390478
$data[11,12] = 33.38796883158906
$data[11,13] = 35.29658989882071
$data[11,14] = 234.8710142410997
$data[11,15] = 61093
$data[11,16] = 303.9452736318408
$data[11,17] = 494.1239844310836
$data[11,18] = 26
$data[11,19] = 6494
$data[11,20] = 2012
$data[11,21] = 2023
$data[12,0] = "aon"
$data[12,1] = "lgbtqiamais"
$data[12,2] = 82
$data[12,3] = 58
$data[12,4] = 0.06142322097378277
$data[12,5] = 0.7073170731707317
$data[12,6] = 1602557.623236534
$data[12,7] = 27630.30384890577
$data[12,8] = 38892.54053994336
$data[12,9] = 721.7894130003107
$data[12,10] = 264456.5194757923
$data[12,11] = 88.86048522662777
$data[12,12] = 39.11867998547866
$data[12,13] = 37.30506273801686
$data[12,14] = 245.6155654729304
$data[12,15] = 17873
$data[12,16] = 308.1551724137931
$data[12,17] = 327.6476927842036
$data[12,18] = 8
$data[12,19] = 1539
$data[12,20] = 2013
$data[12,21] = 2023
$data[13,0] = "aon"
$data[13,1] = "midia_independente"
$data[13,2] = 140
$data[13,3] = 97
$data[13,4] = 0.1048689138576779
$data[13,5] = 0.6928571428571428
$data[13,6] = 4218883.099892105
$data[13,7] = 43493.64020507325
$data[13,8] = 57897.47966668471
$data[13,9] = 1405.425637867093
$data[13,10] = 264585.9073482947
$data[13,11] = 103.6833090125289
$data[13,12] = 51.99423496130147
$data[13,13] = 36.80839302979295
$data[13,14] = 323.2845357010965
$data[13,15] = 38250
$data[13,16] = 394.3298969072165
$data[13,17] = 391.0591046986345
$data[13,18] = 29
$data[13,19] = 1711
$data[13,20] = 2012
$data[13,21] = 2023
$data[14,0] = "aon"
$data[14,1] = "politica"
$data[14,2] = 240
$data[14,3] = 154
$data[14,4] = 0.1797752808988764
$data[14,5] = 0.6416666666666667
$data[14,6] = 4982733.189867401
$data[14,7] = 32355.4103238143
$data[14,8] = 34162.57496563493
$data[14,9] = 54.53892516702949
$data[14,10] = 215281.2939270826
$data[14,11] = 95.62742936226397
$data[14,12] = 44.41979710275996
$data[14,13] = 33.80063482849972
$data[14,14] = 362.0414364166904
$data[14,15] = 51893
$data[14,16] = 336.9675324675325
$data[14,17] = 298.7578952675915
$data[14,18] = 1
$data[14,19] = 1588
$data[14,20] = 2012
$data[14,21] = 2023
$data[15,0] = "aon"
$data[15,1] = "questoes_genero"
$data[15,2] = 35
$data[15,3] = 24
$data[15,4] = 0.02621722846441948
$data[15,5] = 0.6857142857142857
$data[15,6] = 622913.4565716616
$data[15,7] = 25954.72735715257
$data[15,8] = 26184.67661836704
$data[15,9] = 3366.13723259077
$data[15,10] = 123112.7022480959
$data[15,11] = 88.50596984923921
$data[15,12] = 44.08134994103607
$data[15,13] = 46.83761258476419
$data[15,14] = 245.6155654729304
$data[15,15] = 7569
$data[15,16] = 315.375
$data[15,17] = 313.5216257016551
$data[15,18] = 51
$data[15,19] = 1489
$data[15,20] = 2013
$data[15,21] = 2023
$data[16,0] = "aon"
$data[16,1] = "religiosidade"
$data[16,2] = 360
$data[16,3] = 211
$data[16,4] = 0.2696629213483146
$data[16,5] = 0.5861111111111111
$data[16,6] = 6620001.799310843
$data[16,7] = 31374.41611047793
$data[16,8] = 59614.08605582522
$data[16,9] = 322.2027789591561
$data[16,10] = 679297.6600721752
$data[16,11] = 87.04136070513805
$data[16,12] = 35.12732957218113
$data[16,13] = 21.61624650544615
$data[16,14] = 301.8203940790075
$data[16,15] = 75649
$data[16,16] = 358.5260663507109
$data[16,17] = 555.7348485891068
$data[16,18] = 4
$data[16,19] = 6494
$data[16,20] = 2012
$data[16,21] = 2023
$data[17,0] = "aon"
$data[17,1] = "saloes_humor"
$data[17,2] = 16
$data[17,3] = 12
$data[17,4] = 0.01198501872659176
$data[17,5] = 0.75
$data[17,6] = 162569.9259774528
$data[17,7] = 13547.4938314544
$data[17,8] = 10287.95417495333
$data[17,9] = 459.3885785954238
$data[17,10] = 29349.35154822051
$data[17,11] = 75.00500439706998
$data[17,12] = 43.45490371809562
$data[17,13] = 32.2694534583262
$data[17,14] = 185.9579322823807
$data[17,15] = 2243
$data[17,16] = 186.9166666666667
$data[17,17] = 124.1256389406485
$data[17,18] = 11
$data[17,19] = 411
$data[17,20] = 2013
$data[17,21] = 2023
$data[18,0] = "aon"
$data[18,1] = "terror"
$data[18,2] = 496
$data[18,3] = 327
$data[18,4] = 0.3715355805243445
$data[18,5] = 0.6592741935483871
$data[18,6] = 9479119.38973001
$data[18,7] = 28988.13269030584
$data[18,8] = 51783.64819049771
$data[18,9] = 787.1021772339901
$data[18,10] = 679297.6600721752
$data[18,11] = 91.19226916374137
$data[18,12] = 57.41534544136429
$data[18,13] = 21.61624650544615
$data[18,14] = 792.0360759681182
$data[18,15] = 108579
$data[18,16] = 332.045871559633
$data[18,17] = 527.8918493761878
$data[18,18] = 16
$data[18,19] = 6494
$data[18,20] = 2012
$data[18,21] = 2023
$data[19,0] = "aon"
$data[19,1] = "webformatos"
$data[19,2] = 119
$data[19,3] = 89
$data[19,4] = 0.08913857677902622
$data[19,5] = 0.7478991596638656
$data[19,6] = 2104894.58419056
$data[19,7] = 23650.50094596134
$data[19,8] = 22257.37070968975
$data[19,9] = 3458.599440871384
$data[19,10] = 136747.5985390203
$data[19,11] = 83.093358381024
$data[19,12] = 29.94469803232121
$data[19,13] = 21.61624650544615
$data[19,14] = 194.2230576381307
$data[19,15] = 25369
$data[19,16] = 285.0449438202247
$data[19,17] = 251.1921534544801
$data[19,18] = 38
$data[19,19] = 1547
$data[19,20] = 2012
$data[19,21] = 2023
$data[20,0] = "aon"
$data[20,1] = "zine"
$data[20,2] = 164
$data[20,3] = 112
$data[20,4] = 0.1228464419475655
$data[20,5] = 0.6829268292682927
$data[20,6] = 2315707.644801917
$data[20,7] = 20675.96111430283
$data[20,8] = 20218.37616086783
$data[20,9] = 54.53892516702949
$data[20,10] = 161153.6281071717
$data[20,11] = 82.11629402200307
$data[20,12] = 32.47932839143289
$data[20,13] = 32.2694534583262
$data[20,14] = 245.6155654729304
$data[20,15] = 28547
$data[20,16] = 254.8839285714286
$data[20,17] = 226.6887672930878
$data[20,18] = 1
$data[20,19] = 1588
$data[20,20] = 2012
$data[20,21] = 2023
$ws.Range("A2:V22").Value = $data

# ---- 3) Number formats for the data columns ----
# (these reuse the workbook's existing style slots: 1 = "#,##0",
#  2 = "0.00%", 3 = "R$ #,##0.00" -- no new styles are created)
$ws.Range("C2:D22").NumberFormat = "#,##0"
$ws.Range("E2:F22").NumberFormat = "0.00%"
$ws.Range("G2:O22").NumberFormat = "R$ #,##0.00"
$ws.Range("P2:T22").NumberFormat = "#,##0"

